$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Logistic Regression
$ws.Range("B2").Value = 0.8017057569296375
$ws.Range("C2").Value = 0.6557377049180327
$ws.Range("D2").Value = 0.5347593582887701
$ws.Range("E2").Value = 0.8983543078412392
$ws.Range("F2").Value = 0.5891016200294551
$ws.Range("G2").Value = 0.8407783777067986

# Row 3: Decision Tree
$ws.Range("B3").Value = 0.7874911158493249
$ws.Range("C3").Value = 0.5944584382871536
$ws.Range("D3").Value = 0.6310160427807486
$ws.Range("E3").Value = 0.8441432720232332
$ws.Range("F3").Value = 0.6121919584954605
$ws.Range("G3").Value = 0.8140844122564981

# Row 4: Random Forest
$ws.Range("B4").Value = 0.7619047619047619
$ws.Range("C4").Value = 0.5365853658536586
$ws.Range("D4").Value = 0.7647058823529411
$ws.Range("E4").Value = 0.7608906098741529
$ws.Range("F4").Value = 0.6306504961411246
$ws.Range("G4").Value = 0.8361658841130397

# Row 5: Gradient Boosting
$ws.Range("B5").Value = 0.8002842928216063
$ws.Range("C5").Value = 0.6715867158671587
$ws.Range("D5").Value = 0.4866310160427808
$ws.Range("E5").Value = 0.9138431752178122
$ws.Range("F5").Value = 0.5643410852713179
$ws.Range("G5").Value = 0.8376554969431229
